# Rename the embedded picture "file names" recorded on the inline
# pictures that live in the headers/footers of this BTEC assignment
# brief: the two Pearson Edexcel logos (in the default + first-page
# footers) swap from "image2.png" to "image1.png", and the BTEC logo
# (in the first-page header) swaps from "image1.jpg" to "image2.jpg".
#
# Setting InlineShape.Name while the handle is scoped to the whole
# header/footer Range can leave a stale handle on a multi-paragraph
# footer story, so we first find the specific paragraph that owns the
# picture and rename it from that narrower Range instead.

function Rename-InlinePicture($headerFooter, [string]$newName) {
    $owner = $null
    $paragraphs = $headerFooter.Range.Paragraphs
    for ($i = 1; $i -le $paragraphs.Count; $i++) {
        $para = $paragraphs.Item($i)
        if ($para.Range.InlineShapes.Count -gt 0) {
            $owner = $para
        }
    }
    $owner.Range.InlineShapes(1).Name = $newName
}

$d = $word.ActiveDocument

# Default footer ("Prepared By: BTEC Internal Assessment QDAM...") -
# Pearson logo: image2.png -> image1.png
Rename-InlinePicture $d.Sections(1).Footers(1) "image1.png"

# First-page footer ("Prepared By: QDAM BTEC Assessment...") -
# Pearson logo: image2.png -> image1.png
Rename-InlinePicture $d.Sections(1).Footers(2) "image1.png"

# First-page header - BTEC logo: image1.jpg -> image2.jpg
Rename-InlinePicture $d.Sections(1).Headers(2) "image2.jpg"
